$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- New hours logged for the week of 2022-12-15 through 2022-12-23 ---
# Row 89 (2022-12-15)
$ws.Range("G89").Value = 3.5
$ws.Range("H89").Value = 3.5
$ws.Range("J89").Value = 1.5

# Row 90 (2022-12-16)
$ws.Range("G90").Value = 2
$ws.Range("J90").Value = 2.5

# Row 94 (2022-12-20)
$ws.Range("G94").Value = 2
$ws.Range("H94").Value = 2
$ws.Range("I94").Value = 2
$ws.Range("J94").Value = 3

# Row 95 (2022-12-21)
$ws.Range("H95").Value = 2.75
$ws.Range("I95").Value = 1.75
$ws.Range("J95").Value = 4.25

# Row 96 (2022-12-22)
$ws.Range("G96").Value = 2
$ws.Range("H96").Value = 2
$ws.Range("I96").Value = 2
$ws.Range("J96").Value = 2

# Row 97 (2022-12-23)
$ws.Range("G97").Value = 2.5
$ws.Range("J97").Value = 1.5

# Update the frozen-pane scroll position / selection to reflect the newly
# entered data further down the sheet.
$excel.ActiveWindow.FreezePanes = $true
$excel.ActiveWindow.ScrollRow = 82
$ws.Range("K133").Select()

$wb.Save()
